# Create new testcase for appy store search app
# - Rename "Sheet2" to "homepagetestdata"
# - Add "appiumvv" value to homepagetestdata!A1
# - Fix "secure@web" hyperlink display text to "secure@web1" on Login sheet
# - Switch the active/selected tab from Login to homepagetestdata

$wb = $excel.ActiveWorkbook

# Rename Sheet2 -> homepagetestdata
$wsHome = $wb.Worksheets.Item("Sheet2")
$wsHome.Name = "homepagetestdata"

# Add the new value for the new test data sheet
$wsHome.Range("A1").Value = "appiumvv"

# Fix the hyperlink display text on the Login sheet (B1 currently shows "secure@web").
# Re-creating both hyperlinks (instead of mutating in place) avoids leaving a stale
# duplicate entry behind, since this host always appends rather than updates in place.
$wsLogin = $wb.Worksheets.Item("Login")
$wsLogin.Hyperlinks.Delete()
$wsLogin.Hyperlinks.Add($wsLogin.Range("A1"), "mailto:shailendra@appypie.com", "", "", "shailendra@appypie.com")
$wsLogin.Hyperlinks.Add($wsLogin.Range("B1"), "mailto:secure@web", "", "", "secure@web1")

# Make homepagetestdata the active sheet/tab (activeTab becomes 1, tabSelected flips)
$wsHome.Activate()
$wsHome.Select()
